$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 15

$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"

# Date column - mirror the format/style used by the other rows in column D
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 4).Value = (Get-Date -Year 2022 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100101
$ws.Cells.Item($newRow, 8).Value = "Berries"
$ws.Cells.Item($newRow, 9).Value = 100101004
$ws.Cells.Item($newRow, 10).Value = "Frambuesa"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 40
$ws.Cells.Item($newRow, 14).Value = 6000
$ws.Cells.Item($newRow, 15).Value = 6000
$ws.Cells.Item($newRow, 16).Value = 6000
$ws.Cells.Item($newRow, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia de Linares"
$ws.Cells.Item($newRow, 19).Value = 3000
$ws.Cells.Item($newRow, 20).Value = 2
